$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 9).Value = 1494
$ws.Cells.Item(3, 9).Value = 64
$ws.Cells.Item(4, 9).Value = 1879
$ws.Cells.Item(5, 9).Value = 2030
$ws.Cells.Item(6, 9).Value = 1874
$ws.Cells.Item(7, 9).Value = 3153
$ws.Cells.Item(8, 9).Value = 2052
$ws.Cells.Item(9, 9).Value = 3174
$ws.Cells.Item(10, 9).Value = 3150
$ws.Cells.Item(11, 9).Value = 3002
$ws.Cells.Item(12, 9).Value = 290
$ws.Cells.Item(13, 9).Value = 79
$ws.Cells.Item(14, 9).Value = 3438
$ws.Cells.Item(16, 9).Value = 3585
$ws.Cells.Item(17, 9).Value = 3586
$ws.Cells.Item(18, 9).Value = 3586
$ws.Cells.Item(19, 9).Value = 3440
$ws.Cells.Item(20, 9).Value = 3579
$ws.Cells.Item(21, 9).Value = 3583
$ws.Cells.Item(22, 9).Value = 3895
$ws.Cells.Item(23, 9).Value = 3804
$ws.Cells.Item(24, 9).Value = 4191
$ws.Cells.Item(25, 9).Value = 3778
$ws.Cells.Item(26, 9).Value = 827
$ws.Cells.Item(27, 9).Value = 4303
$ws.Cells.Item(28, 9).Value = 773
$ws.Cells.Item(29, 9).Value = 4561
$ws.Cells.Item(30, 9).Value = 4710
$ws.Cells.Item(31, 9).Value = 4713
$ws.Cells.Item(32, 9).Value = 2356
$ws.Cells.Item(33, 9).Value = 106
$ws.Cells.Item(34, 9).Value = 3672
$ws.Cells.Item(35, 9).Value = 106
$ws.Cells.Item(36, 9).Value = 4405
$ws.Cells.Item(37, 9).Value = 5618
$ws.Cells.Item(38, 9).Value = 6410
$ws.Cells.Item(39, 9).Value = 5405
$ws.Cells.Item(40, 9).Value = 4289
$ws.Cells.Item(41, 9).Value = 6411
$ws.Cells.Item(42, 9).Value = 3965
$ws.Cells.Item(43, 9).Value = 1768
$ws.Cells.Item(44, 9).Value = 5076
$ws.Cells.Item(45, 9).Value = 6429
$ws.Cells.Item(47, 9).Value = 913
$ws.Cells.Item(48, 9).Value = 6574
$ws.Cells.Item(49, 9).Value = 6550
$ws.Cells.Item(50, 9).Value = 6572
$ws.Cells.Item(51, 9).Value = 6581
$ws.Cells.Item(52, 9).Value = 5594
$ws.Cells.Item(54, 9).Value = 2021
$ws.Cells.Item(55, 9).Value = 5974
$ws.Cells.Item(56, 9).Value = 1421
$ws.Cells.Item(57, 9).Value = 616
$ws.Cells.Item(58, 9).Value = 1663
$ws.Cells.Item(59, 9).Value = 4566
